# Changing the parameters for locate r_peak function
$wb = $excel.ActiveWorkbook

# Rename the second sheet ("sheet1") to "agsuy"
$ws = $wb.Worksheets.Item(2)
$ws.Name = "agsuy"

# Updated values for column J (rows 1-48)
$jValues = @(
    35.66036176681519,
    35.27510190010071,
    35.63436627388,
    35.85895705223083,
    36.01353073120117,
    36.64841961860657,
    36.25137758255005,
    36.87122178077698,
    35.86889672279358,
    36.32108521461487,
    36.12348079681396,
    36.08053636550903,
    35.7022557258606,
    35.80611824989319,
    35.24312305450439,
    36.49177050590515,
    36.28306865692139,
    37.03895664215088,
    35.94142127037048,
    35.18150424957275,
    36.30312657356262,
    35.24853181838989,
    35.65342211723328,
    36.1482834815979,
    35.22927808761597,
    35.49432039260864,
    36.83135080337524,
    35.55908012390137,
    35.77238821983337,
    36.59109735488892,
    36.44655609130859,
    35.42531585693359,
    37.54795575141907,
    38.58660840988159,
    35.66844010353088,
    36.98784446716309,
    36.35196089744568,
    37.52797031402588,
    36.07941365242004,
    35.14920425415039,
    35.61071491241455,
    35.67176651954651,
    35.98828673362732,
    36.20944476127625,
    35.7984356880188,
    35.46314167976379,
    36.8673210144043,
    35.98003172874451
)

for ($i = 0; $i -lt $jValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
